$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header D1 from "autofill_Sentiment" to "Language"
$ws.Range("D1").Value = "Language"

# Row 2: mark as reviewed (TRUE) and set its Language to French
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = "French"

# Clear out the old autofilled Sentiment values for the remaining rows (3-20)
$ws.Range("D3:D20").ClearContents()
